$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2025-10-07T08:55:00+00:00"

# --- Elements sheet: remove the "valueString" comment slice row ---
$ws = $wb.Worksheets.Item("Elements")

# Row 74 is "QuestionnaireResponse.item.item.answer.value[x]:valueString"
# (the commentaire slice). Deleting it shifts every following row up by one.
$ws.Rows.Item(74).Delete()

# Row 73 ("...answer.value[x]") was the sliced element; now that its only
# slice is gone, clear its slicing discriminator/description/rules.
$ws.Cells.Item(73, 28).Value = ""
$ws.Cells.Item(73, 29).Value = ""
$ws.Cells.Item(73, 31).Value = ""
